$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add a new "2021" data column (R) mirroring the existing "2020" column (Q),
# copying formatting from the corresponding Q-column cell in each row and
# then writing the new value for that row.

$newValues = @{
    4  = 2021
    5  = 99.4
    6  = 98.1
    7  = 99.319469393395053
    8  = 99.442213297634979
    9  = 99.1
    10 = 99.3
    11 = 99.799160124155549
    12 = 99.3
    13 = 99.538370126605429
    14 = 99.765563948945029
}

foreach ($row in 4..14) {
    $srcCell = $ws.Range("Q$row")
    $dstCell = $ws.Range("R$row")

    $srcCell.Copy()
    $dstCell.PasteSpecial(-4122)  # xlPasteFormats
    $dstCell.Value = $newValues[$row]
}

$ws.Application.CutCopyMode = $false

$ws.Range("U4").Select()
